$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new attendance/log row (row 70) below the existing data.
$ws.Range("A70").Value = 200785
$ws.Range("B70").Value = "General Surgery"
$ws.Range("C70").Value = 45906
$ws.Range("D70").Value = 0.460590277777778

# Copy the time-format (h:mm:ss) style from the row above instead of
# fabricating a brand-new number format / style entry.
$ws.Range("D69").Copy()
$ws.Range("D70").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("E70").Value = "Edited"
$ws.Range("F70").Value = "system"

# Match the saved selection state: active cell on the newly written F70.
$ws.Range("F70").Select()
